# 06组项目计划表 - add the "第十三周三" status block (rows 209-218) and
# fill in the completion-rate values for the previous week's rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Append a brand-new weekly block (10 rows) for 第十三周三, reusing the
#    exact same layout/formatting as the previous week's block (rows
#    199-208: date banner, column headers, 5 member rows, blank row, and the
#    2-row "总结" footer) by copying it wholesale, then overwriting the
#    handful of cells that differ. This has to happen BEFORE step 2 fills in
#    the "完成情况" column below, otherwise those values would get dragged
#    along into the freshly copied block.
# ---------------------------------------------------------------------------
$ws.Range("A199:D208").Copy($ws.Range("A209:D218"))

# ---------------------------------------------------------------------------
# 2) Fill in the "完成情况" (completion rate) column for the 第十三周一 block
#    (rows 201-205) that was left blank before.
# ---------------------------------------------------------------------------
$ws.Range("C201").Value = 0.8
$ws.Range("C202").Value = 1
$ws.Range("C203").Value = 1
$ws.Range("C204").Value = 1
$ws.Range("C205").Value = 1

# Row 209: date banner for the new block.
$ws.Range("A209").Value = "日期：2018.11.21 第十三周三"

# Row 211: 李福森's entry keeps his name but gets this week's content.
$ws.Range("B211").Value = "内容:修改不符合美观的页面"

# Rows 212-214: the other three members only have their names carried
# forward for this block - no content/notes recorded yet.
$ws.Range("B212").Value = ""
$ws.Range("B213").Value = ""
$ws.Range("B214").Value = ""

# Row 215 (previously 王一鸣's row when copied from row 205) is blank in the
# new block, same as row 216.
$ws.Range("A215").Value = ""
$ws.Range("B215").Value = ""
